# Weekly injuries check-in (Week #15):
#  - Chad Bettis: refresh his Last.Updated date + injury details (still out,
#    now expected back "later part of July" instead of "around the All-Star
#    break").
#  - Ian Desmond's calf-strain entry is replaced by a brand new entry for
#    Tyler Chatwood, who has just landed on the 10-day DL with a calf strain.
#  - Rows are re-sorted alphabetically by player Name (Anderson, Bettis,
#    Chatwood, Dahl).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Tyler Anderson (unchanged)
$ws.Cells.Item(2, 1).Value = "Tyler Anderson"
$ws.Cells.Item(2, 2).Value = "anderty01"
$ws.Cells.Item(2, 3).Value = "June 27 2017"
$ws.Cells.Item(2, 4).Value = "Knee"
$ws.Cells.Item(2, 5).Value = "Anderson has been designated for the 10-day disabled list as he will have arthroscopic left knee surgery and is expected to miss up to four weeks of action."

# Row 3 - Chad Bettis (updated date + details)
$ws.Cells.Item(3, 1).Value = "Chad Bettis"
$ws.Cells.Item(3, 2).Value = "bettich01"
$ws.Cells.Item(3, 3).Value = "July 10 2017"
$ws.Cells.Item(3, 4).Value = "Illness"
$ws.Cells.Item(3, 5).Value = "Bettis is on the 60-day disabled list while recovering from testicular cancer but is expected to make his season debut sometime during the later part of July."

# Row 4 - Tyler Chatwood (new, replaces Ian Desmond)
$ws.Cells.Item(4, 1).Value = "Tyler Chatwood"
$ws.Cells.Item(4, 2).Value = "chatwty01"
$ws.Cells.Item(4, 3).Value = "July 16 2017"
$ws.Cells.Item(4, 4).Value = "Calf"
$ws.Cells.Item(4, 5).Value = "Chatwood has been placed on the 10-day disabled list with a right calf strain and there is no timetable for return."

# Row 5 - David Dahl (unchanged, moved down one row)
$ws.Cells.Item(5, 1).Value = "David Dahl"
$ws.Cells.Item(5, 2).Value = "dahlda01"
$ws.Cells.Item(5, 3).Value = "June 27 2017"
$ws.Cells.Item(5, 4).Value = "Ribs"
$ws.Cells.Item(5, 5).Value = "Dahl has been moved to the 60-day disabled list with a stress reaction of his sixth rib and it is unknown as to when he will be ready to join the lineup."

# Reflect the updated selection/cursor position left behind by the edit.
[void]$ws.Range("A20").Select()
